# Auto-generated edit script: applies 89 cell updates (text values, percent changes,
# and two pairs of row re-ordering via swapped Coin/Link text) to the crypto table.
# Every write goes through Set-TextValue so cells keep their original "General" style
# (no numberFormat/quotePrefix residue) while still landing as literal text, even when
# the new value is something Excel would otherwise auto-convert to a number (e.g. "1.00").

function Set-TextValue {
    param(
        $ws,
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
Set-TextValue $ws "D2" "69.342.69"
Set-TextValue $ws "E2" "  -2.00%  "

# Row 3
Set-TextValue $ws "D3" "3.680.13"
Set-TextValue $ws "E3" "  -3.48%  "

# Row 4
Set-TextValue $ws "E4" "  +0.01%  "

# Row 5
Set-TextValue $ws "D5" "683.55"
Set-TextValue $ws "E5" "  -3.65%  "

# Row 6
Set-TextValue $ws "D6" "162.50"
Set-TextValue $ws "E6" "  -4.69%  "

# Row 7
Set-TextValue $ws "D7" "3.677.40"

# Row 8
Set-TextValue $ws "E8" "  +0.03%  "

# Row 9
Set-TextValue $ws "E9" "  -4.36%  "

# Row 10
Set-TextValue $ws "E10" "  -7.44%  "

# Row 11
Set-TextValue $ws "D11" "7.21"
Set-TextValue $ws "E11" "  -4.26%  "

# Row 12
Set-TextValue $ws "E12" "  -1.85%  "

# Row 13
Set-TextValue $ws "E13" "  -6.20%  "

# Row 14
Set-TextValue $ws "E14" "  -6.73%  "

# Row 15
Set-TextValue $ws "D15" "4.303.88"
Set-TextValue $ws "E15" "  -3.45%  "

# Row 16
Set-TextValue $ws "D16" "3.679.81"
Set-TextValue $ws "E16" "  -3.22%  "

# Row 17
Set-TextValue $ws "D17" "69.377.52"
Set-TextValue $ws "E17" "  -2.05%  "

# Row 19
Set-TextValue $ws "D19" "16.34"
Set-TextValue $ws "E19" "  -5.43%  "

# Row 20
Set-TextValue $ws "E20" "  -6.74%  "

# Row 21
Set-TextValue $ws "D21" "478.14"
Set-TextValue $ws "E21" "  -3.58%  "

# Row 22
Set-TextValue $ws "D22" "9.83"
Set-TextValue $ws "E22" "  -7.48%  "

# Row 23
Set-TextValue $ws "E23" "  -8.05%  "

# Row 24
Set-TextValue $ws "D24" "80.12"
Set-TextValue $ws "E24" "  -4.87%  "

# Row 25
Set-TextValue $ws "D25" "3.826.17"
Set-TextValue $ws "E25" "  -3.49%  "

# Row 26
Set-TextValue $ws "E26" "  -10.88%  "

# Row 27
Set-TextValue $ws "B27" "InternetComputer(DFINITY)"
Set-TextValue $ws "C27" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D27" "11.52"
Set-TextValue $ws "E27" "  -4.75%  "

# Row 28
Set-TextValue $ws "B28" "Dai"
Set-TextValue $ws "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D28" "1.00"
Set-TextValue $ws "E28" "  +0.02%  "

# Row 29
Set-TextValue $ws "D29" "9.59"
Set-TextValue $ws "E29" "  -7.78%  "

# Row 30
Set-TextValue $ws "E30" "  -10.07%  "

# Row 31
Set-TextValue $ws "E31" "  -10.81%  "

# Row 32
Set-TextValue $ws "E32" "  -5.84%  "

# Row 33
Set-TextValue $ws "D33" "6.88"
Set-TextValue $ws "E33" "  -6.23%  "

# Row 34
Set-TextValue $ws "D34" "27.13"
Set-TextValue $ws "E34" "  -6.84%  "

# Row 35
Set-TextValue $ws "E35" "  +0.19%  "

# Row 36
Set-TextValue $ws "E36" "  -5.03%  "

# Row 37
Set-TextValue $ws "D37" "3.643.41"
Set-TextValue $ws "E37" "  -3.67%  "

# Row 38
Set-TextValue $ws "E38" "  -5.60%  "

# Row 39
Set-TextValue $ws "D39" "6.13"
Set-TextValue $ws "E39" "  +2.95%  "

# Row 40
Set-TextValue $ws "D40" "0.0940"
Set-TextValue $ws "E40" "  -7.39%  "

# Row 41
Set-TextValue $ws "E41" "  +0.01%  "

# Row 42
Set-TextValue $ws "D42" "2.15"
Set-TextValue $ws "E42" "  -6.31%  "

# Row 43
Set-TextValue $ws "E43" "  -0.10%  "

# Row 44
Set-TextValue $ws "D44" "0.957"
Set-TextValue $ws "E44" "  -8.02%  "

# Row 45
Set-TextValue $ws "D45" "48.20"
Set-TextValue $ws "E45" "  -1.28%  "

# Row 46
Set-TextValue $ws "B46" "dogwifhat"
Set-TextValue $ws "C46" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D46" "2.85"
Set-TextValue $ws "E46" "  -11.91%  "

# Row 47
Set-TextValue $ws "B47" "Monero"
Set-TextValue $ws "C47" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D47" "154.12"
Set-TextValue $ws "E47" "  -6.85%  "

# Row 48
Set-TextValue $ws "B48" "FLOKI"
Set-TextValue $ws "C48" "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws "D48" "0.000283"
Set-TextValue $ws "E48" "  -12.34%  "

# Row 49
Set-TextValue $ws "B49" "ONDO"
Set-TextValue $ws "C49" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws "D49" "1.33"
Set-TextValue $ws "E49" "  -1.69%  "

# Row 50
Set-TextValue $ws "D50" "392.76"
Set-TextValue $ws "E50" "  -7.18%  "
